$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "58.931.21"
$ws.Range("E2").Value = "  +0.00%  "
Set-TextValue "D3" "2.508.37"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue "D5" "532.15"
$ws.Range("E5").Value = "  -0.69%  "
Set-TextValue "D6" "135.05"
$ws.Range("E6").Value = "  -2.17%  "
$ws.Range("E7").Value = "  +0.12%  "
Set-TextValue "D8" "0.566"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("E12").Value = "  -0.63%  "
Set-TextValue "D13" "2.953.16"
$ws.Range("E13").Value = "  -1.03%  "
Set-TextValue "D14" "58.826.68"
$ws.Range("E14").Value = "  -0.27%  "
Set-TextValue "D15" "22.69"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("E16").Value = "  -1.34%  "
Set-TextValue "D17" "2.513.14"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("E19").Value = "  -0.65%  "
Set-TextValue "D20" "322.05"
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("E22").Value = "  +0.41%  "
Set-TextValue "D23" "65.05"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("E26").Value = "  -1.10%  "
Set-TextValue "D27" "7.54"
$ws.Range("E27").Value = "  -1.28%  "
Set-TextValue "D28" "0.0₃0761"
$ws.Range("E28").Value = "  -2.25%  "
Set-TextValue "D29" "6.49"
$ws.Range("E29").Value = "  -3.62%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "1.75"
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D31" "169.46"
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  -4.61%  "
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("E37").Value = "  -3.09%  "
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("E39").Value = "  -3.98%  "
Set-TextValue "D40" "280.67"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("E42").Value = "  -0.41%  "
Set-TextValue "D43" "129.78"
$ws.Range("E43").Value = "  -0.27%  "
Set-TextValue "D44" "4.98"
$ws.Range("E44").Value = "  -5.76%  "
Set-TextValue "D46" "0.0925"
$ws.Range("E46").Value = "  -0.77%  "
Set-TextValue "D47" "0.0499"
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("E48").Value = "  -3.04%  "
Set-TextValue "D49" "17.21"
$ws.Range("E49").Value = "  -1.17%  "
Set-TextValue "D50" "1.759.28"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("E51").Value = "  -0.61%  "
